$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the 19th day's "Lön" (payday) row as complete: bump the day count
# and apply the same green highlight fill used by the rest of that row.
$ws.Range("I19").Value = 4
$ws.Range("G19").Interior.Color = $ws.Range("D19").Interior.Color

# Update the active selection to reflect where the user left off editing.
$ws.Range("L19").Select()
